# Daily attendance processing - 2025-10-05 07:37:45
# Applies the session-analysis refresh: two sessions move from
# "Pending" to "Recorded" (rows 8 and 37), their stats ripple through
# the summary / group-statistics tables, and a batch of "Recorded By"
# e-mail lists get re-ordered (same people, new order) across many rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Row 8  (Year2 / A1 / PARASITOLOGY / Session 2) : Pending -> Recorded
#    Copy the look of an already-"Recorded" row (row 7) so the same
#    shared cell style (green fill) is reused instead of creating a
#    brand new style entry.
# ---------------------------------------------------------------------
$ws.Range("A7:I7").Copy()
$ws.Range("A8:I8").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G8").Value = "Alshimaa_khaled@med.asu.edu.eg"
$ws.Range("H8").Value = "133/217"
$ws.Range("I8").Value = "Recorded"

# ---------------------------------------------------------------------
# 2) Row 37 (Year2 / A3 / BIOCHEMISTRY LAB/CBL / Session 1) : Pending -> Recorded
# ---------------------------------------------------------------------
$ws.Range("A38:I38").Copy()
$ws.Range("A37:I37").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G37").Value = "Kerelos.zareef@med.asu.edu.eg"
$ws.Range("H37").Value = "108/220"
$ws.Range("I37").Value = "Recorded"

# ---------------------------------------------------------------------
# 3) Overall dashboard statistics (K/L columns, rows 6-10)
#    The percentage columns are stored as plain TEXT (e.g. "28.8%"),
#    not numeric percentages, in the original workbook. A leading
#    apostrophe forces text-entry (exactly like typing it in Excel),
#    then we re-apply the clean "style 5" number format (General) from
#    a neighboring text cell so no stray style with quotePrefix lingers
#    on the cell itself.
# ---------------------------------------------------------------------
$ws.Range("L6").Value  = 44        # Recorded Sessions
$ws.Range("L8").Value  = 101       # Pending Sessions

$ws.Range("L9").Value  = "'28.8%"  # Coverage %
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("L10").Value = "'48.6%"  # Average Attendance %
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 4) Group statistics table (rows 15 & 17 -> PHARMACOLOGY A1 / PHYSIOLOGY A3)
# ---------------------------------------------------------------------
$ws.Range("O15").Value = 5
$ws.Range("Q15").Value = 11

$ws.Range("R15").Value = "'29.4%"
$ws.Range("M15").Copy()
$ws.Range("R15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("S15").Value = "'64.1%"
$ws.Range("M15").Copy()
$ws.Range("S15").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("S16").Value = "'59.2%"
$ws.Range("N16").Copy()
$ws.Range("S16").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("O17").Value = 6
$ws.Range("Q17").Value = 10

$ws.Range("R17").Value = "'35.3%"
$ws.Range("M17").Copy()
$ws.Range("R17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("S17").Value = "'50.9%"
$ws.Range("M17").Copy()
$ws.Range("S17").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 5) Students column text update (row 25, already "Recorded")
# ---------------------------------------------------------------------
$ws.Range("H25").Value = "124/216"

# ---------------------------------------------------------------------
# 6) "Recorded By" e-mail list re-orderings (same attendees, new order)
# ---------------------------------------------------------------------
$ws.Range("G14").Value  = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G31").Value  = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"

$ws.Range("G17").Value  = "Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"
$ws.Range("G34").Value  = "Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg"

$ws.Range("G35").Value  = "Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"

$ws.Range("G45").Value  = "System, backup@backdoor.com, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G62").Value  = "System, backup@backdoor.com, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

$ws.Range("G51").Value  = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"
$ws.Range("G68").Value  = "Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg"

$ws.Range("G72").Value  = "wessam.atef@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

$ws.Range("G76").Value  = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"

$ws.Range("G83").Value  = "afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G150").Value = "afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg"

$ws.Range("G85").Value  = "ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"
$ws.Range("G102").Value = "ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

$ws.Range("G98").Value  = "user@user.com, afaf.abdallah@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg"

$ws.Range("G99").Value  = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G149").Value = "user@user.com, Walaa.h.ghanima@med.asu.edu.eg"

$ws.Range("G116").Value = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"
$ws.Range("G133").Value = "enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg"

$ws.Range("G119").Value = "aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G136").Value = "aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
